# Adds a processing/creation date stamp to the "About" sheet (cell C1),
# stored as an Excel date serial (44307 = 4/21/2021) formatted with the
# built-in short-date number format (numFmtId 14).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("C1").Value = 44307
$ws.Range("C1").NumberFormat = "mm-dd-yy"
